$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric/percent-looking strings stay as text (matching inlineStr source)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.51"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.50%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.31"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.63%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.159"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.66%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05649"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.98%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.470"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.04%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8180"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.53%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8338"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.18%"
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0005999"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.36%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1330"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.28%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06925"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.26%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02931"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.86%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09398"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.08%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001507"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.34%"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.04280"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-8.98%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006152"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.89%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.509"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.45%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.02%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.311"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "9.10%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.03100"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-3.22%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1292"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.13%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.126"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-43.25%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1374"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.10%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001225"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.00%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-2.78%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.00009798"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "2.02%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.00007257"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-47.79%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03649"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.09%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006069"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "78.23%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1054"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-21.99%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002300"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-13.60%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008197"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.04%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005382"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.62%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.02%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1010"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-23.96%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002655"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "29.17%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.02%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.02%"
